## Insert a new weekly price record as row 536 on the "Ciboulette" sheet,
## pushing the existing rows 536-652 down to 537-653 (dimension grows from
## A1:R652 to A1:R653).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 536..652 down by one (copies formatting/style from the row
# above, matching default Excel "Insert" behaviour - e.g. column D keeps
# its date style).
$ws.Rows.Item(536).Insert()

# Populate the newly inserted row 536 with the new weekly record.
$ws.Range("A536").Value = 9
$ws.Range("B536").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C536").Value = "Metropolitana"
$ws.Range("D536").Value = 45244
$ws.Range("E536").Value = 13
$ws.Range("F536").Value = 100112039
$ws.Range("G536").Value = "Ciboulette"
$ws.Range("H536").Value = "Sin especificar"
$ws.Range("I536").Value = "Primera"
$ws.Range("J536").Value = 430
$ws.Range("K536").Value = 1500
$ws.Range("L536").Value = 1500
$ws.Range("M536").Value = 1500
$ws.Range("N536").Value = "`$/docena de atados"
$ws.Range("O536").Value = "Región Metropolitana"
$ws.Range("P536").Value = 500
$ws.Range("Q536").Value = 3
$ws.Range("R536").Value = "Hortaliza"
